$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44313, 12, 55, 313.0513973476009),
    @(44314, 1, 49, 278.9003358187717),
    @(44315, 2, 43, 244.7492742899425),
    @(44316, 2, 31, 176.4471512322841),
    @(44317, 12, 37, 210.5982127611133),
    @(44318, 5, 41, 233.3655871136661)
)

$startRow = 239
$lastRow = 238
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]

    $ws.Cells.Item($lastRow, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
}
